# Translate short-doc headings/body text from Vietnamese to French.
$d = $word.ActiveDocument

function Find-ParagraphByText {
    # Locate the (single) paragraph whose text - minus its trailing
    # paragraph mark - equals $Text exactly. Searching by content (rather
    # than a hard-coded paragraph index) keeps this robust to any shift in
    # paragraph numbering.
    param([string]$Text)
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Length -eq ($Text.Length + 1) -and $t.Substring(0, $Text.Length) -eq $Text) {
            return $p
        }
    }
    throw "Paragraph with text '$Text' not found"
}

function Replace-InParagraph {
    # Plain Range.Text assignment (NOT Find.Execute) so that straight
    # apostrophes in the French replacement text are preserved verbatim
    # instead of being "smart-quoted" by AutoCorrect-style Find/Replace.
    param(
        [string]$OldText,
        [string]$NewText
    )
    $p = Find-ParagraphByText $OldText
    $full = $p.Range
    $r = $d.Range($full.Start, $full.Start + $OldText.Length)
    $r.Text = $NewText
}

function Split-ParagraphRun {
    # Paragraph currently holds: [Run1: $OldFirst (bold)] [Run2: empty, non-bold]
    # After: [Run1: $NewFirst (bold, same rPr)] [Run2: $NewSecond (non-bold, same rPr)]
    param(
        [string]$OldFirst,
        [string]$NewFirst,
        [string]$NewSecond
    )
    $p = Find-ParagraphByText $OldFirst
    $full = $p.Range
    $paraStart = $full.Start
    $paraEnd = $full.End

    # Insert the new second-run text right before the paragraph mark, then
    # demote it to non-bold so the engine splits it into its own run
    # (matching the pre-existing, non-bold, empty trailing run).
    $insertPos = $paraEnd - 1
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertBefore($NewSecond)
    $r2 = $d.Range($insertPos, $insertPos + $NewSecond.Length)
    $r2.Font.Bold = 0

    # Now replace the first run's original text with its new (shorter) text.
    $r1 = $d.Range($paraStart, $paraStart + $OldFirst.Length)
    $r1.Text = $NewFirst
}

# --- Title page ---
Replace-InParagraph "Kế hoạch Chiến lược" "Plan stratégique"
Replace-InParagraph "Tháng 1 năm 2050" "Janvier 2050"

# --- Table of contents heading ---
Replace-InParagraph "Mục lục" "Table des matières"

# --- Executive summary ---
Replace-InParagraph "Tóm tắt Điều hành" "Résumé exécutif"
Replace-InParagraph `
    "Kế hoạch Chiến lược Đơn giản của Fortune Flow phác thảo tầm nhìn, mục tiêu và chiến lược của chúng tôi trong ba năm tới. Lộ trình ngắn gọn này sẽ hướng dẫn quỹ đạo tăng trưởng của công ty và nâng cao lợi thế cạnh tranh của chúng tôi. Bằng cách ưu tiên đổi mới, tính bền vững và sự hài lòng của khách hàng, chúng tôi sẵn sàng đạt được tiến bộ đáng kể trong sự hiện diện trên thị trường." `
    "Le Plan Stratégique Simple Fortune Flow décrit notre vision, nos objectifs et nos stratégies pour les trois prochaines années. Cette feuille de route concise guidera la trajectoire de croissance de notre entreprise et renforcera notre avantage concurrentiel. En privilégiant l'innovation, la durabilité et la satisfaction client, nous sommes prêts à réaliser des progrès substantiels dans notre présence sur le marché."

# --- Vision and Mission ---
Replace-InParagraph "Tầm nhìn và Sứ mệnh" "Vision et Mission"
Split-ParagraphRun `
    "Tầm nhìn: Fortune Flow khao khát trở thành nhà lãnh đạo toàn cầu trong các giải pháp công nghệ bền vững, thúc đẩy tiến bộ cho cả xã hội và môi trường. Mục tiêu của chúng tôi là cách mạng hóa các ngành công nghiệp thông qua các công nghệ tiên tiến tối ưu hóa việc sử dụng tài nguyên và giảm tác động sinh thái." `
    "Vision :" `
    " Fortune Flow aspire à devenir un leader mondial dans les solutions technologiques durables, stimulant le progrès pour la société et l'environnement. Notre objectif est de révolutionner les industries grâce à des technologies avancées qui optimisent l'utilisation des ressources et réduisent l'impact écologique."

# --- Objectives ---
Replace-InParagraph "Mục tiêu/Các Mục tiêu" "Objectifs"
Replace-InParagraph "Mục tiêu 1: Lãnh đạo Thị trường" "Objectif 1 : Leadership sur le marché"
Replace-InParagraph "Mục tiêu 1.1: Tăng Thị phần" "Objectif 1.1 : Augmenter la part de marché"
Split-ParagraphRun `
    "Mục tiêu: Đạt được mức tăng 15% thị phần trong vòng ba năm." `
    "Objectif : " `
    "Atteindre une augmentation de 15 % de la part de marché en trois ans."

Replace-InParagraph "Mục tiêu 1.2: Mở rộng Sự hiện diện tại Các Thị trường Mới nổi" "Objectif 1.2 : Étendre la présence sur les marchés émergents"
Split-ParagraphRun `
    "Mục tiêu: Mở rộng sự hiện diện của chúng tôi tại các thị trường mới nổi, nhắm đến mức tăng trưởng 20% doanh thu từ các khu vực mới." `
    "Objectif : Ét" `
    "endre notre présence sur les marchés émergents, avec une croissance de 20 % du chiffre d'affaires provenant de nouvelles régions."

Write-Output "done"
